$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for the 8c254ff2 (row 3) de-de entry
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-20 00:52:42"

# Sheet "zh-cn": row 3 (8c254ff2 file) Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-20 00:52:37"
$wsZhCn.Range("K3").Value = "2016-08-20 00:52:53"

# Sheet "de-de": row 3 (8c254ff2 file) Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-20 00:52:42"
$wsDeDe.Range("K3").Value = "2016-08-20 00:52:59"
